$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("heavy")

# Update the input values that drive the loadtest numbers (corrected per
# Confluence "Performance Test Data.xls"). Downstream formulas (I5, R2, S2,
# U2, X2, F6, I6, J6, F7, I7, ...) recalculate automatically.
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("M6").Value = 0

$ws.Range("G7").Value = 0.25
$ws.Range("M7").Value = 0.25

# Restore the active selection on the "heavy" sheet to I11.
$ws.Activate()
[void]$ws.Range("I11").Select()
